$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated average_county_temperature (column I) values with NOAA data
$ws.Range("I2").Value  = 19.30324074074072
$ws.Range("I3").Value  = 17.25771604938272
$ws.Range("I4").Value  = 13.62268518518517
$ws.Range("I11").Value = 13.75752314814816
$ws.Range("I12").Value = 19.79629629629628
$ws.Range("I13").Value = 0.8611111111111096
$ws.Range("I15").Value = 0.8611111111111096
$ws.Range("I17").Value = 5.486111111111112
$ws.Range("I18").Value = 16.86342592592595
$ws.Range("I19").Value = 16.86342592592595
$ws.Range("I20").Value = 5.486111111111112
$ws.Range("I21").Value = 12.41429539295394
$ws.Range("I23").Value = 19.36574074074073
$ws.Range("I29").Value = 12.41429539295394
$ws.Range("I31").Value = 12.41429539295394
$ws.Range("I32").Value = 12.41429539295394
$ws.Range("I33").Value = 12.41429539295394

# Downstream recalculated values that depend on column I (worst/best ASHP COP)
$ws.Range("N32").Value = 1.466311090415359
$ws.Range("O32").Value = 1.558048815385048
